$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.139.76"
$ws.Range("E2").Value = "  +3.45%  "
$ws.Range("D3").Value = "2.363.41"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.87"
$ws.Range("E6").Value = "  -1.15%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  -0.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.72"
$ws.Range("E10").Value = "  -0.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0916"
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.43"
$ws.Range("E12").Value = "  -2.06%  "
$ws.Range("E13").Value = "  +1.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.974"
$ws.Range("E14").Value = "  -2.80%  "
$ws.Range("D15").Value = "2.724.37"
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.21"
$ws.Range("E16").Value = "  -1.99%  "
$ws.Range("D17").Value = "2.362.74"
$ws.Range("E17").Value = "  +1.42%  "
$ws.Range("D18").Value = "45.110.92"
$ws.Range("E18").Value = "  +3.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.59"
$ws.Range("E19").Value = "  +10.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.26"
$ws.Range("E20").Value = "  -4.37%  "
$ws.Range("E21").Value = "  -1.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.11"
$ws.Range("E22").Value = "  -1.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.48"
$ws.Range("E23").Value = "  -0.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "258.95"
$ws.Range("E24").Value = "  -3.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.29"
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.08"
$ws.Range("E27").Value = "  -0.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.22"
$ws.Range("E28").Value = "  -5.59%  "
$ws.Range("E29").Value = "  +2.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0968"
$ws.Range("E30").Value = "  +8.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.34"
$ws.Range("E31").Value = "  -1.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "37.28"
$ws.Range("E32").Value = "  -3.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "168.97"
$ws.Range("E33").Value = "  +0.88%  "
$ws.Range("E34").Value = "  +3.88%  "
$ws.Range("E35").Value = "  -1.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.116"
$ws.Range("E36").Value = "  +3.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.71"
$ws.Range("E37").Value = "  -1.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.93"
$ws.Range("E38").Value = "  +3.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.89"
$ws.Range("E39").Value = "  +1.78%  "
$ws.Range("E40").Value = "  -3.64%  "
$ws.Range("E41").Value = "  +1.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.49"
$ws.Range("E42").Value = "  -5.18%  "
$ws.Range("E43").Value = "  -3.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "69.13"
$ws.Range("E44").Value = "  -3.80%  "
$ws.Range("E45").Value = "  -3.30%  "
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("D47").Value = "1.827.97"
$ws.Range("E47").Value = "  +9.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "81.14"
$ws.Range("E48").Value = "  +5.35%  "
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.54"
$ws.Range("E49").Value = "  +3.61%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "111.77"
$ws.Range("E50").Value = "  -2.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.13"
$ws.Range("E51").Value = "  +1.99%  "
